# Generate Report for Handoff
#
# The localization run for 4d1e0d08-ad0f-401d-a993-5f2f03cd48cf.md has
# progressed: it is now "Ready for handoff" (was "Handed back: in sync with
# en-US"), with fresh handoff timestamps, and an Error Detail noting the
# handback file in use is not the latest version available upstream.
# 3800af99-368a-41d2-a59a-8bf899ade174.md (the other row) is unaffected.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/73b5d6cd04f96a6fd40bf4377efc543c943d2eab/e2e/4d1e0d08-ad0f-401d-a993-5f2f03cd48cf.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/774ac177e81b3916e0d515c984a8d07516229e46/e2e/4d1e0d08-ad0f-401d-a993-5f2f03cd48cf.md."

# --- Overview sheet: row 3 is the 4d1e0d08 file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-01 10:55:53"

# --- zh-cn sheet: row 3 is the 4d1e0d08 file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-09-01 10:55:49"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Range("P1").ColumnWidth = 39.17

# --- de-de sheet: row 3 is the 4d1e0d08 file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-09-01 10:55:53"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Range("P1").ColumnWidth = 39.17
